# Generate Report for Handoff
# Adds a new row (row 3) describing the handoff of
# "b9cccdbb-31bc-4522-bc8d-ac4d8d4f85ea.md" to each of the report sheets:
# "Overview", "zh-cn" and "de-de". Each sheet backs an Excel Table
# (ListObject), so we grow the table by one row, fill in the new data,
# restore the date/time number format on the datetime columns and add the
# hyperlink to the source file - mirroring how the existing row (for
# 603796dc-fce3-494f-9d55-e04697052c8e.md) is laid out.

$wb = $excel.ActiveWorkbook

$newFile        = "b9cccdbb-31bc-4522-bc8d-ac4d8d4f85ea.md"
$newPath        = "e2e\b9cccdbb-31bc-4522-bc8d-ac4d8d4f85ea.md"
$newFileUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f16fdbaa46a1b729355cb025932a4faa4d43c8c2/e2e/b9cccdbb-31bc-4522-bc8d-ac4d8d4f85ea.md"
$zhXlf          = "b9cccdbb-31bc-4522-bc8d-ac4d8d4f85ea.d4ce27c15a4b28fa77217b99a11a3e92daaf84a5.zh-cn.xlf"
$deXlf          = "b9cccdbb-31bc-4522-bc8d-ac4d8d4f85ea.d4ce27c15a4b28fa77217b99a11a3e92daaf84a5.de-de.xlf"
$hoGenerateDate = "2016-08-17 04:35:55"
$zhHandoffDate  = "2016-08-17 04:35:50"
$deHandoffDate  = "2016-08-17 04:35:55"
$dateFormat     = "yyyy-mm-dd HH:mm:ss"
$epoch          = "0001-01-01 00:00:00"

# Helper: write a literal empty string into a cell (instead of leaving it
# completely blank) so the cell participates in the row like the original
# template rows do ("" entries still get a shared-string cell).
function Set-EmptyText($range) {
    $range.Value = "'"
    $range.Style = "Normal"
}

# Helper: write True/False as literal text (not as a native boolean) -
# matches how the template stores these as shared strings.
function Set-BoolText($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------
# Overview sheet: File Name | Path And Name | Extension | Publish URL |
#                 zh-cn | de-de | Latest HO Xliff Generate Date
# ---------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$loOv = $wsOv.ListObjects.Item(1)
$loOv.ListRows.Add() | Out-Null

$wsOv.Range("A3").Value = $newFile
$wsOv.Range("B3").Value = $newPath
$wsOv.Range("C3").Value = ".md"
Set-EmptyText $wsOv.Range("D3")
$wsOv.Range("E3").Value = "Ready for handoff"
$wsOv.Range("F3").Value = "Ready for handoff"
$wsOv.Range("G3").Value = $hoGenerateDate
$wsOv.Range("G3").NumberFormat = $dateFormat

$wsOv.Hyperlinks.Add($wsOv.Range("B3"), $newFileUrl, "", "", $newPath) | Out-Null
$wsOv.Range("B3").Style = "Hyperlink"

# ---------------------------------------------------------------
# Per-locale sheets ("zh-cn" and "de-de"):
# Source File Name | File Extension | Status | Source Path | Priority |
# Content Duplicate | Latest Handoff File | Latest Handoff Datetime |
# Latest Target File | Latest Handback File | Latest Handback DateTime |
# Reference Tokens | To be localized | Dependency From | Has metadata |
# Error Detail
# ---------------------------------------------------------------
function Add-LocaleRow($sheetName, $handoffXlf, $handoffDate) {
    $ws = $wb.Worksheets.Item($sheetName)
    $lo = $ws.ListObjects.Item(1)
    $lo.ListRows.Add() | Out-Null

    $ws.Range("A3").Value = $newFile
    $ws.Range("B3").Value = ".md"
    $ws.Range("C3").Value = "Ready for handoff"
    $ws.Range("D3").Value = "e2e"
    $ws.Range("E3").Value = "ht"
    Set-BoolText $ws.Range("F3") "False"
    $ws.Range("G3").Value = $handoffXlf
    $ws.Range("H3").Value = $handoffDate
    $ws.Range("H3").NumberFormat = $dateFormat
    Set-EmptyText $ws.Range("I3")
    Set-EmptyText $ws.Range("J3")
    $ws.Range("K3").Value = $epoch
    $ws.Range("K3").NumberFormat = $dateFormat
    Set-EmptyText $ws.Range("L3")
    Set-BoolText $ws.Range("M3") "True"
    Set-EmptyText $ws.Range("N3")
    Set-BoolText $ws.Range("O3") "False"
    Set-EmptyText $ws.Range("P3")

    $ws.Hyperlinks.Add($ws.Range("A3"), $newFileUrl, "", "", $newFile) | Out-Null
    $ws.Range("A3").Style = "Hyperlink"
}

Add-LocaleRow "zh-cn" $zhXlf $zhHandoffDate
Add-LocaleRow "de-de" $deXlf $deHandoffDate
